# Update attendance/price figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 2215
    $ws.Range("G2").Value = 70

    $ws.Range("G3").Value = 60

    $ws.Range("F4").Value = 1607

    $ws.Range("F5").Value = 7454
}
